$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (co2)
$ws.Range("B2").Value = 340542
$ws.Range("C2").Value = 534.0941587234468
$ws.Range("D2").Value = 116.5276101818023
$ws.Range("F2").Value = 449
$ws.Range("G2").Value = 496
$ws.Range("H2").Value = 584

# Row 3 (humidity)
$ws.Range("B3").Value = 340542
$ws.Range("C3").Value = 44.90025967428392
$ws.Range("D3").Value = 4.62925788268382
$ws.Range("F3").Value = 41.45
$ws.Range("G3").Value = 44.84
$ws.Range("H3").Value = 48.12

# Row 4 (pm25)
$ws.Range("B4").Value = 340542
$ws.Range("C4").Value = 1.705674747901874
$ws.Range("D4").Value = 3.238081598277234
$ws.Range("F4").Value = 0.61
$ws.Range("G4").Value = 1.33
$ws.Range("H4").Value = 2.33

# Row 5 (pressure)
$ws.Range("B5").Value = 340542
$ws.Range("C5").Value = 322.5149398312103
$ws.Range("D5").Value = 9.867247193481946
$ws.Range("F5").Value = 318.08
$ws.Range("G5").Value = 324.39
$ws.Range("H5").Value = 329.84

# Row 6 (temperature)
$ws.Range("B6").Value = 340542
$ws.Range("C6").Value = 22.51477858825049
$ws.Range("D6").Value = 1.835523249557064
$ws.Range("F6").Value = 21.52
$ws.Range("G6").Value = 22.19
$ws.Range("H6").Value = 23.13

# Row 7 (rssi)
$ws.Range("B7").Value = 340542
$ws.Range("C7").Value = -76.33234079790452
$ws.Range("D7").Value = 23.56758322221495

# Row 8 (snr)
$ws.Range("B8").Value = 339939
$ws.Range("C8").Value = 7.720823147682379
$ws.Range("D8").Value = 6.579145716284592

# Row 9 (SF)
$ws.Range("B9").Value = 340542
$ws.Range("C9").Value = 9.324559085222969
$ws.Range("D9").Value = 1.688312520378005

# Row 10 (frequency)
$ws.Range("B10").Value = 340542
$ws.Range("C10").Value = 867.8304931550293
$ws.Range("D10").Value = 0.4611783823689606

# Row 11 (f_count)
$ws.Range("B11").Value = 340530
$ws.Range("C11").Value = 15176.54487416674
$ws.Range("D11").Value = 11053.5498440258
$ws.Range("F11").Value = 5787
$ws.Range("G11").Value = 13032
$ws.Range("H11").Value = 24291

# Row 12 (p_count)
$ws.Range("B12").Value = 340542
$ws.Range("C12").Value = 16969.28617321799
$ws.Range("D12").Value = 12312.16610049215
$ws.Range("F12").Value = 6532
$ws.Range("G12").Value = 14595
$ws.Range("H12").Value = 27149

# Row 13 (toa)
$ws.Range("B13").Value = 340542
$ws.Range("C13").Value = 0.5574502936142973
$ws.Range("D13").Value = 0.5908322098315822

# Row 14 (distance)
$ws.Range("B14").Value = 340542
$ws.Range("C14").Value = 23.88875381010272
$ws.Range("D14").Value = 13.41276591323907

# Row 15 (c_walls)
$ws.Range("B15").Value = 340542
$ws.Range("C15").Value = 0.6732121148052223
$ws.Range("D15").Value = 0.7486168012273926

# Row 16 (w_walls)
$ws.Range("B16").Value = 340542
$ws.Range("C16").Value = 1.829788983444039
$ws.Range("D16").Value = 1.669302736950451

# Row 17 (exp_pl)
$ws.Range("B17").Value = 340542
$ws.Range("C17").Value = 93.73234079790458
$ws.Range("D17").Value = 23.56758322220557

# Row 18 (n_power)
$ws.Range("B18").Value = 339939
$ws.Range("C18").Value = -85.51379225515126
$ws.Range("D18").Value = 21.33041486398002
$ws.Range("F18").Value = -102.4139268515822
$ws.Range("G18").Value = -85.02214159641585

# Row 19 (esp)
$ws.Range("B19").Value = 339939
$ws.Range("C19").Value = -77.79296910746885
$ws.Range("D19").Value = 25.63942277324374
$ws.Range("F19").Value = -92.66683163887967
$ws.Range("G19").Value = -74.26572375596102
$ws.Range("H19").Value = -55.2778545523916
